$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L1").Select()
$excel.ActiveWindow.Split = $true
Write-Host "SplitCol: " $excel.ActiveWindow.SplitColumn
